$d = $word.ActiveDocument

# Locate the paragraph that holds the italic "2 Kings" subtitle that
# immediately follows the "2KI" Heading2 abbreviation paragraph, and
# remove it (paragraph mark included) so the heading is followed
# directly by the next paragraph (the single-space one).
$count = $d.Paragraphs.Count
$targetIndex = -1

for ($i = 1; $i -lt $count; $i++) {
    $cur = $d.Paragraphs.Item($i)
    $curText = $cur.Range.Text.Trim()
    if ($curText -eq "2KI") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        $nextText = $nextPara.Range.Text.Trim()
        if ($nextText -eq "2 Kings") {
            $targetIndex = $i + 1
            break
        }
    }
}

if ($targetIndex -gt 0) {
    $p = $d.Paragraphs.Item($targetIndex)
    $p.Range.Delete()
    Write-Host "Removed '2 Kings' paragraph (index $targetIndex) following the '2KI' heading."
} else {
    Write-Host "Target paragraph not found; no changes made."
}
